$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date in column C for rows 2-11 from 45224 to 45233
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = 45233
}
